$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-29 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-09-30 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("44-39=", $true, $true, $false, $false, $false, $true, 1, $false, "24+59=", 2) | Out-Null
$d.Content.Find.Execute("63-53=", $true, $true, $false, $false, $false, $true, 1, $false, "61+10=", 2) | Out-Null
$d.Content.Find.Execute("41-2=", $true, $true, $false, $false, $false, $true, 1, $false, "88-67=", 2) | Out-Null
$d.Content.Find.Execute("99-39=", $true, $true, $false, $false, $false, $true, 1, $false, "26-12=", 2) | Out-Null
$d.Content.Find.Execute("2+18=", $true, $true, $false, $false, $false, $true, 1, $false, "70+6=", 2) | Out-Null
$d.Content.Find.Execute("56-24=", $true, $true, $false, $false, $false, $true, 1, $false, "71-3=", 2) | Out-Null
$d.Content.Find.Execute("57-46=", $true, $true, $false, $false, $false, $true, 1, $false, "29-26=", 2) | Out-Null
$d.Content.Find.Execute("64+31=", $true, $true, $false, $false, $false, $true, 1, $false, "63-56=", 2) | Out-Null
$d.Content.Find.Execute("1+28=", $true, $true, $false, $false, $false, $true, 1, $false, "24+42=", 2) | Out-Null
$d.Content.Find.Execute("45+49=", $true, $true, $false, $false, $false, $true, 1, $false, "81-50=", 2) | Out-Null
$d.Content.Find.Execute("66-3=", $true, $true, $false, $false, $false, $true, 1, $false, "77+12=", 2) | Out-Null
$d.Content.Find.Execute("36+39=", $true, $true, $false, $false, $false, $true, 1, $false, "78-38=", 2) | Out-Null
$d.Content.Find.Execute("27-19=", $true, $true, $false, $false, $false, $true, 1, $false, "63+16=", 2) | Out-Null
$d.Content.Find.Execute("11+72=", $true, $true, $false, $false, $false, $true, 1, $false, "61+19=", 2) | Out-Null
$d.Content.Find.Execute("40-14=", $true, $true, $false, $false, $false, $true, 1, $false, "28-11=", 2) | Out-Null
$d.Content.Find.Execute("48-13=", $true, $true, $false, $false, $false, $true, 1, $false, "41+46=", 2) | Out-Null
$d.Content.Find.Execute("73+18=", $true, $true, $false, $false, $false, $true, 1, $false, "12+60=", 2) | Out-Null
$d.Content.Find.Execute("21+15=", $true, $true, $false, $false, $false, $true, 1, $false, "86-14=", 2) | Out-Null
$d.Content.Find.Execute("92-70=", $true, $true, $false, $false, $false, $true, 1, $false, "38+28=", 2) | Out-Null
$d.Content.Find.Execute("35+7=", $true, $true, $false, $false, $false, $true, 1, $false, "45-7=", 2) | Out-Null
$d.Content.Find.Execute("47+28=", $true, $true, $false, $false, $false, $true, 1, $false, "73+1=", 2) | Out-Null
$d.Content.Find.Execute("43-29=", $true, $true, $false, $false, $false, $true, 1, $false, "79-33=", 2) | Out-Null
$d.Content.Find.Execute("91-49=", $true, $true, $false, $false, $false, $true, 1, $false, "40-39=", 2) | Out-Null
$d.Content.Find.Execute("29+68=", $true, $true, $false, $false, $false, $true, 1, $false, "16+3=", 2) | Out-Null
$d.Content.Find.Execute("80-49=", $true, $true, $false, $false, $false, $true, 1, $false, "18-2=", 2) | Out-Null
$d.Content.Find.Execute("28+27=", $true, $true, $false, $false, $false, $true, 1, $false, "77+1=", 2) | Out-Null
$d.Content.Find.Execute("75-54=", $true, $true, $false, $false, $false, $true, 1, $false, "56+29=", 2) | Out-Null
$d.Content.Find.Execute("6+91=", $true, $true, $false, $false, $false, $true, 1, $false, "52-0=", 2) | Out-Null
$d.Content.Find.Execute("6+64=", $true, $true, $false, $false, $false, $true, 1, $false, "4+86=", 2) | Out-Null
$d.Content.Find.Execute("91-54=", $true, $true, $false, $false, $false, $true, 1, $false, "77-61=", 2) | Out-Null
$d.Content.Find.Execute("62-1=", $true, $true, $false, $false, $false, $true, 1, $false, "92-91=", 2) | Out-Null
$d.Content.Find.Execute("44+39=", $true, $true, $false, $false, $false, $true, 1, $false, "30+12=", 2) | Out-Null
$d.Content.Find.Execute("29+69=", $true, $true, $false, $false, $false, $true, 1, $false, "75-33=", 2) | Out-Null
$d.Content.Find.Execute("9+65=", $true, $true, $false, $false, $false, $true, 1, $false, "47-14=", 2) | Out-Null
$d.Content.Find.Execute("73-59=", $true, $true, $false, $false, $false, $true, 1, $false, "55+10=", 2) | Out-Null
$d.Content.Find.Execute("73-10=", $true, $true, $false, $false, $false, $true, 1, $false, "17-3=", 2) | Out-Null
$d.Content.Find.Execute("51-50=", $true, $true, $false, $false, $false, $true, 1, $false, "80-39=", 2) | Out-Null
$d.Content.Find.Execute("70-21=", $true, $true, $false, $false, $false, $true, 1, $false, "8+40=", 2) | Out-Null
$d.Content.Find.Execute("48+34=", $true, $true, $false, $false, $false, $true, 1, $false, "86-16=", 2) | Out-Null
$d.Content.Find.Execute("96-93=", $true, $true, $false, $false, $false, $true, 1, $false, "49-19=", 2) | Out-Null
$d.Content.Find.Execute("13+41=", $true, $true, $false, $false, $false, $true, 1, $false, "88-57=", 2) | Out-Null
$d.Content.Find.Execute("48-28=", $true, $true, $false, $false, $false, $true, 1, $false, "45+16=", 2) | Out-Null
$d.Content.Find.Execute("12+2=", $true, $true, $false, $false, $false, $true, 1, $false, "88-12=", 2) | Out-Null
$d.Content.Find.Execute("0+98=", $true, $true, $false, $false, $false, $true, 1, $false, "18+38=", 2) | Out-Null
$d.Content.Find.Execute("75-22=", $true, $true, $false, $false, $false, $true, 1, $false, "0+38=", 2) | Out-Null
$d.Content.Find.Execute("86-74=", $true, $true, $false, $false, $false, $true, 1, $false, "41+44=", 2) | Out-Null
$d.Content.Find.Execute("87-27=", $true, $true, $false, $false, $false, $true, 1, $false, "32+26=", 2) | Out-Null
$d.Content.Find.Execute("27+3=", $true, $true, $false, $false, $false, $true, 1, $false, "23+68=", 2) | Out-Null
$d.Content.Find.Execute("83-43=", $true, $true, $false, $false, $false, $true, 1, $false, "75+7=", 2) | Out-Null
$d.Content.Find.Execute("27+29=", $true, $true, $false, $false, $false, $true, 1, $false, "63+31=", 2) | Out-Null
$d.Content.Find.Execute("15+21=", $true, $true, $false, $false, $false, $true, 1, $false, "32+16=", 2) | Out-Null
$d.Content.Find.Execute("34+2=", $true, $true, $false, $false, $false, $true, 1, $false, "58-4=", 2) | Out-Null
$d.Content.Find.Execute("34+55=", $true, $true, $false, $false, $false, $true, 1, $false, "95-62=", 2) | Out-Null
$d.Content.Find.Execute("51-48=", $true, $true, $false, $false, $false, $true, 1, $false, "17+41=", 2) | Out-Null
$d.Content.Find.Execute("77-14=", $true, $true, $false, $false, $false, $true, 1, $false, "38+20=", 2) | Out-Null
$d.Content.Find.Execute("55+28=", $true, $true, $false, $false, $false, $true, 1, $false, "48-24=", 2) | Out-Null
$d.Content.Find.Execute("57+32=", $true, $true, $false, $false, $false, $true, 1, $false, "47-10=", 2) | Out-Null
$d.Content.Find.Execute("17+25=", $true, $true, $false, $false, $false, $true, 1, $false, "90-15=", 2) | Out-Null
$d.Content.Find.Execute("52+1=", $true, $true, $false, $false, $false, $true, 1, $false, "60+1=", 2) | Out-Null
$d.Content.Find.Execute("2+38=", $true, $true, $false, $false, $false, $true, 1, $false, "79+8=", 2) | Out-Null
$d.Content.Find.Execute("8+32=", $true, $true, $false, $false, $false, $true, 1, $false, "89-61=", 2) | Out-Null
$d.Content.Find.Execute("58-26=", $true, $true, $false, $false, $false, $true, 1, $false, "26+67=", 2) | Out-Null
$d.Content.Find.Execute("18+69=", $true, $true, $false, $false, $false, $true, 1, $false, "37-18=", 2) | Out-Null
$d.Content.Find.Execute("2+45=", $true, $true, $false, $false, $false, $true, 1, $false, "49+0=", 2) | Out-Null
$d.Content.Find.Execute("80-34=", $true, $true, $false, $false, $false, $true, 1, $false, "19+63=", 2) | Out-Null
$d.Content.Find.Execute("40+49=", $true, $true, $false, $false, $false, $true, 1, $false, "10+37=", 2) | Out-Null
$d.Content.Find.Execute("53-41=", $true, $true, $false, $false, $false, $true, 1, $false, "54+22=", 2) | Out-Null
$d.Content.Find.Execute("46+49=", $true, $true, $false, $false, $false, $true, 1, $false, "90-22=", 2) | Out-Null
$d.Content.Find.Execute("12+56=", $true, $true, $false, $false, $false, $true, 1, $false, "91+8=", 2) | Out-Null
$d.Content.Find.Execute("20-11=", $true, $true, $false, $false, $false, $true, 1, $false, "47-16=", 2) | Out-Null
$d.Content.Find.Execute("25-17=", $true, $true, $false, $false, $false, $true, 1, $false, "47+51=", 2) | Out-Null
$d.Content.Find.Execute("13+32=", $true, $true, $false, $false, $false, $true, 1, $false, "45+37=", 2) | Out-Null
$d.Content.Find.Execute("89-41=", $true, $true, $false, $false, $false, $true, 1, $false, "66-9=", 2) | Out-Null
$d.Content.Find.Execute("83+15=", $true, $true, $false, $false, $false, $true, 1, $false, "10+63=", 2) | Out-Null
$d.Content.Find.Execute("45-19=", $true, $true, $false, $false, $false, $true, 1, $false, "21-12=", 2) | Out-Null
$d.Content.Find.Execute("47+36=", $true, $true, $false, $false, $false, $true, 1, $false, "27+30=", 2) | Out-Null
$d.Content.Find.Execute("51+4=", $true, $true, $false, $false, $false, $true, 1, $false, "61-22=", 2) | Out-Null
$d.Content.Find.Execute("10+47=", $true, $true, $false, $false, $false, $true, 1, $false, "20+14=", 2) | Out-Null
$d.Content.Find.Execute("25+23=", $true, $true, $false, $false, $false, $true, 1, $false, "6+79=", 2) | Out-Null
$d.Content.Find.Execute("8+76=", $true, $true, $false, $false, $false, $true, 1, $false, "91-19=", 2) | Out-Null
$d.Content.Find.Execute("46+6=", $true, $true, $false, $false, $false, $true, 1, $false, "99-52=", 2) | Out-Null
$d.Content.Find.Execute("12+36=", $true, $true, $false, $false, $false, $true, 1, $false, "27+67=", 2) | Out-Null
$d.Content.Find.Execute("1+54=", $true, $true, $false, $false, $false, $true, 1, $false, "36+46=", 2) | Out-Null
$d.Content.Find.Execute("55+26=", $true, $true, $false, $false, $false, $true, 1, $false, "21+53=", 2) | Out-Null
$d.Content.Find.Execute("15+54=", $true, $true, $false, $false, $false, $true, 1, $false, "68-65=", 2) | Out-Null
$d.Content.Find.Execute("21+3=", $true, $true, $false, $false, $false, $true, 1, $false, "63-16=", 2) | Out-Null
$d.Content.Find.Execute("63+34=", $true, $true, $false, $false, $false, $true, 1, $false, "3+63=", 2) | Out-Null
$d.Content.Find.Execute("95-11=", $true, $true, $false, $false, $false, $true, 1, $false, "69-22=", 2) | Out-Null
$d.Content.Find.Execute("56-16=", $true, $true, $false, $false, $false, $true, 1, $false, "44-40=", 2) | Out-Null
$d.Content.Find.Execute("26-8=", $true, $true, $false, $false, $false, $true, 1, $false, "80+10=", 2) | Out-Null
$d.Content.Find.Execute("73-39=", $true, $true, $false, $false, $false, $true, 1, $false, "46-0=", 2) | Out-Null
$d.Content.Find.Execute("98-63=", $true, $true, $false, $false, $false, $true, 1, $false, "69+19=", 2) | Out-Null
$d.Content.Find.Execute("89-51=", $true, $true, $false, $false, $false, $true, 1, $false, "11+42=", 2) | Out-Null
$d.Content.Find.Execute("3-2=", $true, $true, $false, $false, $false, $true, 1, $false, "1+31=", 2) | Out-Null
$d.Content.Find.Execute("56+9=", $true, $true, $false, $false, $false, $true, 1, $false, "43-10=", 2) | Out-Null
$d.Content.Find.Execute("15+34=", $true, $true, $false, $false, $false, $true, 1, $false, "11+62=", 2) | Out-Null
$d.Content.Find.Execute("98-62=", $true, $true, $false, $false, $false, $true, 1, $false, "66+23=", 2) | Out-Null
$d.Content.Find.Execute("10+40=", $true, $true, $false, $false, $false, $true, 1, $false, "92-0=", 2) | Out-Null
$d.Content.Find.Execute("40-13=", $true, $true, $false, $false, $false, $true, 1, $false, "30-11=", 2) | Out-Null
$d.Content.Find.Execute("34+40=", $true, $true, $false, $false, $false, $true, 1, $false, "37+41=", 2) | Out-Null
